$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "quantity" column (F) for the future production rows,
# switching from produced qty to uncon_planned_qty while keeping
# today's produced quantity as-is (row 6 unchanged).
$ws.Range("F2").Value = -242
$ws.Range("F3").Value = -546
$ws.Range("F4").Value = -851
$ws.Range("F5").Value = -107
$ws.Range("F7").Value = -107
